# Commit: "Change names from *img to img*"
# Rename the *img sheets so the "img" token moves from a prefix to a
# suffix: himg -> imgh, timg -> imgt, simg -> imgs, gimg -> imgg,
# wimg -> imgw, bimg -> imgb, eimg -> imge.
# This also moves the active/selected tab from "holiday" to the renamed
# "imge" sheet (was "eimg"), matching the workbook's new activeTab/
# tabSelected state.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# Make the renamed "imge" sheet (last tab, was "eimg") the active tab,
# moving the selection away from "holiday".
$wb.Worksheets.Item("imge").Activate()
